$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "93.747.90"
$ws.Range("E2").Value = "  +1.52%  "

# Row 3
$ws.Range("D3").Value = "3.088.48"
$ws.Range("E3").Value = "  -0.73%  "

# Row 4
$ws.Range("E4").Value = "  +0.04%  "

# Row 5
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "233.63"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -3.38%  "

# Row 6
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "610.06"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.85%  "

# Row 7
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "1.09"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +0.07%  "

# Row 8
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.378"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -5.33%  "

# Row 9
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -0.01%  "

# Row 10
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.821"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +11.86%  "

# Row 11
$ws.Range("D11").Value = "3.083.43"
$ws.Range("E11").Value = "  -0.78%  "

# Row 12
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.196"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -3.43%  "

# Row 13
$ws.Range("D13").Value = "93.843.05"
$ws.Range("E13").Value = "  +1.83%  "

# Row 14
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "0.0000239"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -6.53%  "

# Row 15
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "33.97"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -1.59%  "

# Row 16
$ws.Range("D16").Value = "3.672.01"
$ws.Range("E16").Value = "  -0.52%  "

# Row 17
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "5.22"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -5.63%  "

# Row 18
$ws.Range("D18").Value = "3.117.21"
$ws.Range("E18").Value = "  +1.29%  "

# Row 19
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "3.63"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -1.16%  "

# Row 20
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "14.56"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -1.61%  "

# Row 21
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "5.77"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -0.89%  "

# Row 22
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "440.24"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -1.99%  "

# Row 23
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "8.84"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -6.22%  "

# Row 24
$ws.Range("B24").Value = "PEPE"
$ws.Range("C24").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "0.0000192"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -5.45%  "

# Row 25
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "8.26"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +4.79%  "

# Row 26
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "5.52"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -1.54%  "

# Row 27
$ws.Range("B27").Value = "Litecoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "84.82"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -2.51%  "

# Row 28
$ws.Range("B28").Value = "Aptos"
$ws.Range("C28").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "11.90"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +0.91%  "

# Row 29
$ws.Range("D29").Value = "3.270.00"
$ws.Range("E29").Value = "  +0.17%  "

# Row 30
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +0.00%  "

# Row 31
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.254"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +8.19%  "

# Row 32
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.179"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +6.23%  "

# Row 33
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.123"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -11.33%  "

# Row 34
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "9.18"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -1.51%  "

# Row 35
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +0.25%  "

# Row 36
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "7.79"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -3.63%  "

# Row 37
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.158"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -4.48%  "

# Row 38
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "25.52"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -2.67%  "

# Row 39
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "1.88"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -2.18%  "

# Row 40
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.446"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +1.18%  "

# Row 41
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "23.91"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +3.78%  "

# Row 42
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "1.27"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -3.12%  "

# Row 43
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "462.87"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -4.09%  "

# Row 44
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "3.70"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -13.92%  "

# Row 45
$ws.Range("E45").Value = "  +0.00%  "

# Row 46
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "3.09"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -11.39%  "

# Row 47
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "161.19"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +0.99%  "

# Row 48
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.678"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -2.25%  "

# Row 49
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "1.84"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -4.02%  "

# Row 50
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "43.69"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -1.12%  "

# Row 51
$ws.Range("B51").Value = "FirstDigitalUSD"
$ws.Range("C51").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.998"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +0.05%  "
